# Fill in the student's answers in column C of worksheet "4" (sheet4.xml)
# so that the grading formulas in column D / G recompute, matching the
# exam attempt recorded in the upstream commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("4")

$ws.Cells.Item(2,3).Value = 'USBDumper'
$ws.Cells.Item(3,3).Value = 'CAST-128'
$ws.Cells.Item(4,3).Value = 'Key stretching'
$ws.Cells.Item(5,3).Value = 'He can use two internal commands VRFY and EXPN, which provide information about valid users, email addresses, etc'
$ws.Cells.Item(6,3).Value = 'Identify assets and create a baseline -> Vulnerability scan -> Risk assessment -> Remediation -> Verification -> Monitor'
$ws.Cells.Item(7,3).Value = 'Trident'
$ws.Cells.Item(8,3).Value = 'Dictionary attack'
$ws.Cells.Item(11,3).Value = '. (Period)'
$ws.Cells.Item(12,3).Value = 'Hootsuite'
$ws.Cells.Item(13,3).Value = 'sV'
$ws.Cells.Item(14,3).Value = 'filetype:'
$ws.Cells.Item(15,3).Value = 'sA'
$ws.Cells.Item(16,3).Value = 'Passive assessment'
$ws.Cells.Item(18,3).Value = 'Cloudborne attack'
$ws.Cells.Item(19,3).Value = 'https://www.guardster.com'
$ws.Cells.Item(20,3).Value = 'Wardriving'
$ws.Cells.Item(21,3).Value = 'Key escrow'
$ws.Cells.Item(23,3).Value = 'Angela''s public key.'
$ws.Cells.Item(24,3).Value = 'A pop-up appears to a user stating, "Your computer may have been infected with spyware. Click here to install an anti-spyware tool to resolve this issue."'
$ws.Cells.Item(25,3).Value = 'WS-Security'
$ws.Cells.Item(27,3).Value = 'Banner grabbing'
$ws.Cells.Item(28,3).Value = 'MAC flooding'
$ws.Cells.Item(29,3).Value = 'RIPE NCC'
$ws.Cells.Item(30,3).Value = 'Performing content enumeration using a wordlist'
$ws.Cells.Item(31,3).Value = 'Burp Suite'
$ws.Cells.Item(32,3).Value = 'Using wget to perform banner grabbing on the webserver'
$ws.Cells.Item(33,3).Value = 'Union SQL injection'
$ws.Cells.Item(34,3).Value = 'Evilginx'
$ws.Cells.Item(35,3).Value = 'Unspecified proxy activities'
$ws.Cells.Item(36,3).Value = 'Actions on Objective'
$ws.Cells.Item(38,3).Value = 'Credential enumerator'
$ws.Cells.Item(39,3).Value = 'WS-Address spoofing'
$ws.Cells.Item(40,3).Value = 'Power/clock/reset glitching'
$ws.Cells.Item(41,3).Value = '[related:]'
$ws.Cells.Item(43,3).Value = 'Linux OS'
$ws.Cells.Item(44,3).Value = 'During a cyberattack, a hacker corrupts the event logs on all machines'
$ws.Cells.Item(45,3).Value = 'Censys'
$ws.Cells.Item(46,3).Value = 'Operational threat intelligence'
$ws.Cells.Item(48,3).Value = 'Serpent'
$ws.Cells.Item(49,3).Value = '.stm'
$ws.Cells.Item(51,3).Value = 'Buffer Overflow'
$ws.Cells.Item(52,3).Value = 'Internal monologue attack'
$ws.Cells.Item(53,3).Value = 'Spear-phishing attack'
$ws.Cells.Item(54,3).Value = 'Kube-scheduler'
$ws.Cells.Item(55,3).Value = 'Spoofed session flood attack'
$ws.Cells.Item(57,3).Value = 'DNS cache snooping'
$ws.Cells.Item(58,3).Value = 'Downgrade security attack'
$ws.Cells.Item(59,3).Value = 'Watering hole'
$ws.Cells.Item(60,3).Value = 'Whitelist validation'
$ws.Cells.Item(61,3).Value = 'File Transfer Protocol'
$ws.Cells.Item(62,3).Value = 'getsystem'
$ws.Cells.Item(63,3).Value = 'BetterCAP'
$ws.Cells.Item(64,3).Value = 'Detecting the presence of Honeyd honeypots'
$ws.Cells.Item(65,3).Value = 'Verbose failure messages'
$ws.Cells.Item(66,3).Value = 'CeWL'
$ws.Cells.Item(67,3).Value = 'Gray hat'
$ws.Cells.Item(68,3).Value = 'She should check her ARP table and see if there is one IP address with two different MAC addresses'
$ws.Cells.Item(69,3).Value = 'btlejack -f 0x9c68fd30 -t -m 0x1fffffffff'
$ws.Cells.Item(70,3).Value = 'Syhunt Hybrid'
$ws.Cells.Item(71,3).Value = 'STP attack'
$ws.Cells.Item(72,3).Value = 'Log4J'
$ws.Cells.Item(74,3).Value = 'Docker'
$ws.Cells.Item(75,3).Value = 'D'
$ws.Cells.Item(76,3).Value = 'sV'
$ws.Cells.Item(78,3).Value = 'VPN'
$ws.Cells.Item(79,3).Value = 'Searching database statements at the IP address given'
$ws.Cells.Item(80,3).Value = 'Reverse engineering'
$ws.Cells.Item(81,3).Value = 'Scanning networks'
$ws.Cells.Item(82,3).Value = 'Worm'
$ws.Cells.Item(84,3).Value = 'Weaponization'
$ws.Cells.Item(85,3).Value = 'Reverse image search'
$ws.Cells.Item(86,3).Value = 'Pharming'
$ws.Cells.Item(87,3).Value = 'Advanced SMS phishing'
$ws.Cells.Item(88,3).Value = 'IaaS'
$ws.Cells.Item(89,3).Value = 'Netcat'
$ws.Cells.Item(90,3).Value = 'IoTSeeker'
$ws.Cells.Item(91,3).Value = 'Lock-in'
$ws.Cells.Item(92,3).Value = 'PCI-DSS'
$ws.Cells.Item(93,3).Value = 'T'
$ws.Cells.Item(94,3).Value = 'Kernel-level rootkit'
$ws.Cells.Item(95,3).Value = 'External assessment'
$ws.Cells.Item(96,3).Value = 'Evil-twin attack'
$ws.Cells.Item(97,3).Value = 'Yagi antenna'
$ws.Cells.Item(98,3).Value = 'wash'
$ws.Cells.Item(99,3).Value = 'Zero trust network'
$ws.Cells.Item(100,3).Value = 'Out-of-band SQLi'
$ws.Cells.Item(101,3).Value = 'OSINT framework'

# Mirror the author's final on-screen state: sheet "4" becomes the active
# (selected) tab, scrolled so row 94 is at the top, with G105 selected.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 94
$ws.Range("G105").Select()
